$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F (shifts old F:... to G:...)
$ws.Columns.Item(6).Insert()

# Set header for new Address column
$ws.Cells.Item(2, 6).Value = 'Address'

# Populate address values for rows 3-40 (rows 6 and 19 intentionally left blank)
$ws.Cells.Item(3, 6).Value = 'Sangameshwar Comp. Jr College (H S) AmingadHunagund'
$ws.Cells.Item(4, 6).Value = 'S G High School GolageriSindagi'
$ws.Cells.Item(5, 6).Value = 'G H S Avathi'
$ws.Cells.Item(7, 6).Value = 'Hebsur Secondary School HebsurHubli'
$ws.Cells.Item(8, 6).Value = 'G H S BommalapuraKoppa'
$ws.Cells.Item(9, 6).Value = 'G H S Janwada'
$ws.Cells.Item(10, 6).Value = 'S J S High School BhairidevarakoppaHubballi'
$ws.Cells.Item(11, 6).Value = 'Anjuman High School AlmelSindgi'
$ws.Cells.Item(12, 6).Value = 'G H S AngadiMudigere'
$ws.Cells.Item(13, 6).Value = 'G G H S HirisaveChannarayapatna'
$ws.Cells.Item(14, 6).Value = 'Govt. High School AskiSindagi'
$ws.Cells.Item(15, 6).Value = 'Sangameshwar High SchoolKudalsangamHungund'
$ws.Cells.Item(16, 6).Value = 'H M S M Urdu High School'
$ws.Cells.Item(17, 6).Value = 'Sri Adichunchanagiri P U CollegeChannarayapatna'
$ws.Cells.Item(18, 6).Value = 'N F H S Mangalpet'
$ws.Cells.Item(20, 6).Value = 'G H S MattihalliHarapanahalli'
$ws.Cells.Item(21, 6).Value = 'S S P U College (H S) KolharBasavana Bagewadi'
$ws.Cells.Item(22, 6).Value = 'Sri Amrutha Lingeshwara High School M ShivaraChannarayapatna'
$ws.Cells.Item(23, 6).Value = 'G H S ManchanabeleMagadi'
$ws.Cells.Item(24, 6).Value = 'G J C KalasaMudigere'
$ws.Cells.Item(25, 6).Value = 'Gangamma S MargolG H S ShahabadChittapur'
$ws.Cells.Item(26, 6).Value = 'Govt. P B P U College for boysJamkhandi'
$ws.Cells.Item(27, 6).Value = 'G H S MalliJewargi'
$ws.Cells.Item(28, 6).Value = 'S V M High School IlkalHunagund'
$ws.Cells.Item(29, 6).Value = 'G H S HerurKoppa'
$ws.Cells.Item(30, 6).Value = 'G G J C (High School Section) Mudigere'
$ws.Cells.Item(31, 6).Value = 'Govt. J C Hagare Belur'
$ws.Cells.Item(32, 6).Value = 'T M A E S High SchoolNeelagundHarapanahalli'
$ws.Cells.Item(33, 6).Value = 'Govt. P B P U CollegeJamkhandi'
$ws.Cells.Item(34, 6).Value = 'S B High School KalakeriSindagi'
$ws.Cells.Item(35, 6).Value = 'Sri Anjaneya High School NagenahalliKadur'
$ws.Cells.Item(36, 6).Value = 'Warriors High SchoolChallakere'
$ws.Cells.Item(37, 6).Value = 'Govt. High SchoolKaknalBhalki'
$ws.Cells.Item(38, 6).Value = 'G H S BelagihallyChannarayapatna'
$ws.Cells.Item(39, 6).Value = 'G H S BachenahattiMagadi'
$ws.Cells.Item(40, 6).Value = 'G G H S And CompositJunior CollegeJamkhandi'
